$d = $word.ActiveDocument

$replacements = @(
    @("71×12=852", "68×94=6392"),
    @("39×36=1404", "36×36=1296"),
    @("55×15=825", "14×15=210"),
    @("20×53=1060", "93×53=4929"),
    @("14×18=252", "62×25=1550"),
    @("56×60=3360", "37×60=2220"),
    @("79×11=869", "58×79=4582"),
    @("56×91=5096", "95×34=3230"),
    @("33×26=858", "12×32=384"),
    @("76×49=3724", "87×17=1479"),
    @("68×50=3400", "81×74=5994"),
    @("30×50=1500", "95×59=5605"),
    @("28×69=1932", "21×78=1638"),
    @("97×54=5238", "33×40=1320"),
    @("14×47=658", "83×37=3071"),
    @("84×15=1260", "33×25=825"),
    @("42×71=2982", "94×19=1786"),
    @("93×23=2139", "56×82=4592"),
    @("46×96=4416", "95×86=8170"),
    @("51×43=2193", "73×48=3504"),
    @("16×91=1456", "71×48=3408"),
    @("99×70=6930", "74×94=6956"),
    @("53×57=3021", "26×39=1014"),
    @("90×58=5220", "54×31=1674"),
    @("91×20=1820", "31×52=1612")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
